# Insert a new record row above row 37 (shifting existing rows 37-103 down to 38-104)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(37).Insert()

$ws.Cells.Item(37, 1).Value2  = 11
$ws.Cells.Item(37, 2).Value   = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value   = "Bíobío"
$ws.Cells.Item(37, 4).Value2  = 44580
$ws.Cells.Item(37, 5).Value2  = 8
$ws.Cells.Item(37, 6).Value   = "Fruta"
$ws.Cells.Item(37, 7).Value2  = 100103
$ws.Cells.Item(37, 8).Value   = "Frutos de hueso (carozo)"
$ws.Cells.Item(37, 9).Value2  = 100103001
$ws.Cells.Item(37, 10).Value  = "Cereza"
$ws.Cells.Item(37, 11).Value  = "Lapins"
$ws.Cells.Item(37, 12).Value  = "Primera"
$ws.Cells.Item(37, 13).Value2 = 220
$ws.Cells.Item(37, 14).Value2 = 4500
$ws.Cells.Item(37, 15).Value2 = 5000
$ws.Cells.Item(37, 16).Value2 = 4773
$ws.Cells.Item(37, 17).Value  = "$/bandeja 10 kilos"
$ws.Cells.Item(37, 18).Value  = "Provincia de Curicó"
$ws.Cells.Item(37, 19).Value2 = 477
$ws.Cells.Item(37, 20).Value2 = 10
